$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.787.60"
$ws.Range("E2").Value = "  -1.20%  "
$ws.Range("D3").Value = "3.386.21"
$ws.Range("E3").Value = "  -1.90%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "'568.93"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Value = "'141.10"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.69%  "
$ws.Range("D8").Value = "3.385.36"
$ws.Range("E8").Value = "  -1.96%  "
$ws.Range("E9").Value = "  -0.59%  "
$ws.Range("D10").Value = "'7.48"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.20%  "
$ws.Range("E11").Value = "  -1.51%  "
$ws.Range("D12").Value = "'0.395"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.71%  "
$ws.Range("D13").Value = "3.964.36"
$ws.Range("E13").Value = "  -1.91%  "
$ws.Range("D14").Value = "'28.48"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.36%  "
$ws.Range("E15").Value = "  +2.22%  "
$ws.Range("D16").Value = "'0.0000171"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.54%  "
$ws.Range("D17").Value = "3.389.65"
$ws.Range("E17").Value = "  -2.07%  "
$ws.Range("D18").Value = "60.857.35"
$ws.Range("E18").Value = "  -1.35%  "
$ws.Range("E19").Value = "  -1.17%  "
$ws.Range("D20").Value = "'13.97"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.96%  "
$ws.Range("E21").Value = "  -5.81%  "
$ws.Range("D22").Value = "'383.49"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.01%  "
$ws.Range("E23").Value = "  -1.07%  "
$ws.Range("D24").Value = "'73.78"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.69%  "
$ws.Range("D25").Value = "'1.00"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.38%  "
$ws.Range("D26").Value = "'0.0000117"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -5.25%  "
$ws.Range("D27").Value = "3.521.44"
$ws.Range("E27").Value = "  -1.96%  "
$ws.Range("D28").Value = "'0.179"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.06%  "
$ws.Range("D29").Value = "'1.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.16%  "
$ws.Range("D30").Value = "'7.40"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.97%  "
$ws.Range("D31").Value = "'7.98"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.93%  "
$ws.Range("E32").Value = "  -2.86%  "
$ws.Range("E33").Value = "  -1.91%  "
$ws.Range("D35").Value = "'23.58"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.94%  "
$ws.Range("D36").Value = "'6.97"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.59%  "
$ws.Range("D37").Value = "'165.55"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.99%  "
$ws.Range("D38").Value = "3.416.49"
$ws.Range("E38").Value = "  -1.83%  "
$ws.Range("D39").Value = "'4.97"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.12%  "
$ws.Range("E40").Value = "  -4.63%  "
$ws.Range("D41").Value = "'28.25"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.17%  "
$ws.Range("E42").Value = "  -1.47%  "
$ws.Range("D43").Value = "'0.999"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.14%  "
$ws.Range("E44").Value = "  -3.01%  "
$ws.Range("E45").Value = "  -0.84%  "
$ws.Range("D46").Value = "'4.41"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.85%  "
$ws.Range("E47").Value = "  -3.89%  "
$ws.Range("D48").Value = "'1.12"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.99%  "
$ws.Range("D49").Value = "2.486.39"
$ws.Range("E49").Value = "  -4.21%  "
$ws.Range("D50").Value = "'23.53"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.48%  "
$ws.Range("D51").Value = "'6.80"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.64%  "
